$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" header in H1, matching the formatting of the existing
# header row (bold, bordered, centered) by copying the format from G1.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# Fill in the new "Save" column values for rows 2-12.
$saveValues = @(0, 0, 0, 0, 1, 0, 0, 0, 0, 1, 0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
